$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.055.39"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.817.76"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "310.73"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.5009"
$ws.Range("E7").Value = "  -2.35%  "
$ws.Range("D8").Value = "0.3913"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.09916"
$ws.Range("E9").Value = "  +26.60%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "40.79"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "6.431"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").Value = "20.58"
$ws.Range("D14").Value = "1.001"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "1.817.58"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "7.289"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "0.00001137"
$ws.Range("E17").Value = "  +5.80%  "
$ws.Range("D18").Value = "92.35"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "0.06648"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "5.955"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "28.114.21"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "159.02"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.033.47"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.70"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "2.421"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "127.09"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "1.036"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "5.572"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "3.618"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "0.06711"
$ws.Range("E35").Value = "  -5.10%  "
$ws.Range("D36").Value = "0.02341"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "8.925"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D39").Value = "4.961"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "0.6192"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "1.176"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "1.285"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "3.691"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "124.42"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "1.941"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.180"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "0.06788"
$ws.Range("E51").Value = "  -0.49%  "
